$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values in row 8 (CO2 utilization) ---
$ws.Range("C8").Value = 3863706.88078501
$ws.Range("D8").Value = 5627447.51890974
$ws.Range("F8").Value = 3562609.531625094
$ws.Range("G8").Value = 4753632.740854284

# --- Update existing values in row 9 (Bio-based feedstock) ---
$ws.Range("B9").Value = 8586550.000000002
$ws.Range("C9").Value = 3927509.509801868
$ws.Range("D9").Value = 3785952.686250791
$ws.Range("E9").Value = 12882248.87532836
$ws.Range("F9").Value = 4845489.232850321
$ws.Range("G9").Value = 6450022.13154179

# --- Add new rows 13-16 with the same formatting as row 12 ---
$ws.Range("A12:G12").Copy()
$ws.Range("A13:G16").PasteSpecial(-4122)

$ws.Range("A13").Value = "Electrification + Bio-based feedstock"
$ws.Range("B13:G13").Value = 0

$ws.Range("A14").Value = "Conventional + Bio-based feedstock"
$ws.Range("B14:G14").Value = 0

$ws.Range("A15").Value = "Conventional + Bio-based feedstock with CC"
$ws.Range("B15:G15").Value = 0

$ws.Range("A16").Value = "Electrification + Bio-based"
$ws.Range("B16:G16").Value = 0
